# Apply the edit described by the commit "adding lists code + exp":
# adds a new "domain" column (C) to Sheet1 with "animals" for every data
# row, and switches the active sheet/selection from Sheet2 back to Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: add column C ("domain" header + "animals" values) ---
$ws1.Range("C1").Value = "domain"

for ($r = 2; $r -le 47; $r++) {
    $ws1.Cells.Item($r, 3).Value = "animals"
}

# --- Selection / active sheet changes ---
# Sheet1 becomes the active/selected sheet with B4 selected.
$ws1.Activate() | Out-Null
$ws1.Range("B4").Select() | Out-Null

# Sheet2 keeps its own selection (C8) but is no longer the active tab.
$ws2.Range("C8").Select() | Out-Null
$ws1.Activate() | Out-Null
